$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 60.733334
$ws.Range("I8").Value = 60.733334
$ws.Range("K8").Value = 182.200002
$ws.Range("M8").Value = -43.20000199999998

# Row 132
$ws.Range("H132").Value = 11911495
$ws.Range("I132").Value = 13163100
$ws.Range("J132").Value = 21250
$ws.Range("K132").Value = 39489300
$ws.Range("L132").Value = 63750
$ws.Range("M132").Value = -39486770
$ws.Range("N132").Value = -68810

# Row 137
$ws.Range("H137").Value = 9260264
$ws.Range("J137").Value = 1523.6428
$ws.Range("L137").Value = 4570.928400000001
$ws.Range("N137").Value = -9670.928400000001

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4275134.5
$ws.Range("I61").Value = 5051925.5
$ws.Range("J61").Value = 2783.3333
$ws.Range("K61").Value = 5051925.5
$ws.Range("L61").Value = 2783.3333
$ws.Range("M61").Value = -5051713.5
$ws.Range("N61").Value = -3207.3333

# Row 74
$ws.Range("H74").Value = 1537.6897
$ws.Range("I74").Value = 1483.72
$ws.Range("J74").Value = 1875
$ws.Range("K74").Value = 1483.72
$ws.Range("L74").Value = 1875
$ws.Range("M74").Value = -609.72
$ws.Range("N74").Value = -3623

# Row 77
$ws.Range("H77").Value = 1537.6897
$ws.Range("I77").Value = 1483.72
$ws.Range("J77").Value = 1875
$ws.Range("K77").Value = 7418.6
$ws.Range("L77").Value = 9375
$ws.Range("M77").Value = -3050.6
$ws.Range("N77").Value = -18111

# Row 132
$ws.Range("H132").Value = 998239.7
$ws.Range("I132").Value = 1020.5319
$ws.Range("K132").Value = 3061.5957
$ws.Range("M132").Value = -531.5956999999999

# Row 136
$ws.Range("H136").Value = 4275134.5
$ws.Range("I136").Value = 5051925.5
$ws.Range("J136").Value = 2783.3333
$ws.Range("K136").Value = 15155776.5
$ws.Range("L136").Value = 8349.999899999999
$ws.Range("M136").Value = -15153226.5
$ws.Range("N136").Value = -13449.9999

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2781607.2
$ws.Range("I134").Value = 1280.0312
$ws.Range("J134").Value = 13902916
$ws.Range("K134").Value = 3840.0936
$ws.Range("L134").Value = 41708748
$ws.Range("M134").Value = -1305.0936
$ws.Range("N134").Value = -41713818

$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 549.5
$ws.Range("I10").Value = 549.5
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 549.5
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -410.5
$ws.Range("N10").ClearContents()

# Row 31
$ws.Range("H31").Value = 1114219.8
$ws.Range("I31").Value = 1427894.6
$ws.Range("K31").Value = 1427894.6
$ws.Range("M31").Value = -1427599.6

# Row 34
$ws.Range("H34").Value = 1114219.8
$ws.Range("I34").Value = 1427894.6
$ws.Range("K34").Value = 1427894.6
$ws.Range("M34").Value = -1427692.6

# Row 53
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

# Row 58
$ws.Range("H58").Value = 38462364
$ws.Range("I58").Value = 71429560
$ws.Range("J58").Value = 638.6667
$ws.Range("K58").Value = 71429560
$ws.Range("L58").Value = 638.6667
$ws.Range("M58").Value = -71429357
$ws.Range("N58").Value = -1044.6667

# Row 122
$ws.Range("H122").Value = 35715060
$ws.Range("I122").Value = 62500400
$ws.Range("J122").Value = 1275.3334
$ws.Range("K122").Value = 187501200
$ws.Range("L122").Value = 3826.0002
$ws.Range("M122").Value = -187498750
$ws.Range("N122").Value = -8726.0002

# Row 132
$ws.Range("H132").Value = 8548429
$ws.Range("I132").Value = 1128.0952
$ws.Range("K132").Value = 3384.2856
$ws.Range("M132").Value = -854.2856000000002

# Row 134
$ws.Range("H134").Value = 27779056
$ws.Range("I134").Value = 1338.5454
$ws.Range("K134").Value = 4015.6362
$ws.Range("M134").Value = -1480.6362

# Row 136
$ws.Range("H136").Value = 38462364
$ws.Range("I136").Value = 71429560
$ws.Range("J136").Value = 638.6667
$ws.Range("K136").Value = 214288680
$ws.Range("L136").Value = 1916.0001
$ws.Range("M136").Value = -214286130
$ws.Range("N136").Value = -7016.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 890.15
$ws.Range("I131").Value = 830
$ws.Range("J131").Value = 890.75757
$ws.Range("K131").Value = 2490
$ws.Range("L131").Value = 2672.27271
$ws.Range("M131").Value = 2550
$ws.Range("N131").Value = -12752.27271

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 296.69565
$ws.Range("I107").Value = 185
$ws.Range("J107").Value = 418.54544
$ws.Range("K107").Value = 185
$ws.Range("L107").Value = 418.54544
$ws.Range("M107").Value = 1735
$ws.Range("N107").Value = -4258.54544

# Row 132
$ws.Range("H132").Value = 3684.6035
$ws.Range("I132").Value = 1696.0416
$ws.Range("J132").Value = 13229.7
$ws.Range("K132").Value = 5088.1248
$ws.Range("L132").Value = 39689.10000000001
$ws.Range("M132").Value = -2558.1248
$ws.Range("N132").Value = -44749.10000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 108
$ws.Range("H108").Value = 25944
$ws.Range("J108").Value = 25944
$ws.Range("L108").Value = 25944
$ws.Range("N108").Value = -33624

# Row 132
$ws.Range("H132").Value = 26323626
$ws.Range("I132").Value = 37038916
$ws.Range("J132").Value = 22454.818
$ws.Range("K132").Value = 111116748
$ws.Range("L132").Value = 67364.454
$ws.Range("M132").Value = -111114218
$ws.Range("N132").Value = -72424.454

# Row 136
$ws.Range("H136").Value = 61509980
$ws.Range("I136").Value = 31748586
$ws.Range("J136").Value = 111112300
$ws.Range("K136").Value = 95245758
$ws.Range("L136").Value = 333336900
$ws.Range("M136").Value = -95243208
$ws.Range("N136").Value = -333342000

$ws = $wb.Worksheets.Item("WVR")
# Row 13
$ws.Range("H13").Value = 200
$ws.Range("I13").Value = 200
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 200
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -60
$ws.Range("N13").ClearContents()

# Row 104
$ws.Range("H104").Value = 12553.75
$ws.Range("J104").Value = 12553.75
$ws.Range("L104").Value = 12553.75
$ws.Range("N104").Value = -19541.75

# Row 132
$ws.Range("H132").Value = 9827428
$ws.Range("I132").Value = 34483.066
$ws.Range("K132").Value = 103449.198
$ws.Range("M132").Value = -100919.198

# Row 136
$ws.Range("H136").Value = 25003190
$ws.Range("I136").Value = 83336250
$ws.Range("J136").Value = 3305.3572
$ws.Range("K136").Value = 250008750
$ws.Range("L136").Value = 9916.071599999999
$ws.Range("M136").Value = -250006200
$ws.Range("N136").Value = -15016.0716
